$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '54.337.20'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -7.73%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.874.04'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -10.51%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.19%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '474.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -11.21%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '126.22'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.65%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.12%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.873.31'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -10.58%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.403'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -12.06%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.66'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -11.86%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0963'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -15.65%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.330'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -15.90%  '

# Row 13
$ws.Range('E13').Value = '  -4.67%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.362.84'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -10.64%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.98'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -11.04%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '54.225.58'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -7.95%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.865.16'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -10.83%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000134'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -14.76%  '

# Row 19
$ws.Range('E19').Value = '  -10.44%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.50'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -13.69%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.09'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -13.68%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '295.93'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -18.09%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.19%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.444'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -14.32%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '58.82'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -16.17%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.33%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.154'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -10.58%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.06%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0811'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -16.09%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.28'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -11.69%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.59%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.20'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -12.32%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.02'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -12.59%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.62'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -16.11%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.24'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -13.53%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '134.68'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -16.74%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.42'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -14.88%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.21'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -15.29%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '22.94'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -12.53%  '

# Row 40
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0615'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -13.10%  '

# Row 41
$ws.Range('B41').Value = 'RenzoRestakedETH'
$ws.Range('C41').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.896.56'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -10.65%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.996'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.50%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '35.21'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -13.84%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.964'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -12.39%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.600'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -16.03%  '

# Row 46
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.31'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -12.29%  '

# Row 47
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.40'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -15.60%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.042.50'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -11.37%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.32'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -15.31%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -13.92%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0212'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -11.86%  '
